# Accounts.xlsx: add a "Loai the" (card type) column and a couple of new
# account rows.
#
# Existing layout:
#   A1 = "Ten tai khoan"   B1 = "So du"
#   A2 = "kikaho"          B2 = "0"
#
# Target layout:
#   A1 Ten tai khoan   B1 So du   C1 Loai the
#   A2 kikaho          B2 0       C2 Tien mat
#   A3 abcd            B3 15174   C3 The ngan hang
#   A4 abcd3           B4 15174   C4 Tien mat
#   A5 cxz             B5 423     C5 The ngan hang

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$text)
    # Assigning a numeric-looking string (e.g. "15174") straight to .Value
    # coerces it to a real number, which is not what the target sheet wants
    # -- every cell there is a plain text value. A leading apostrophe forces
    # Excel to store it as text instead ("number stored as text"), then
    # ClearFormats() drops the resulting quote-prefix number format so the
    # cell keeps the sheet's normal default style.
    $range.Value = "'" + $text
    $range.ClearFormats()
}

# New header cell for the third column.
$ws.Range("C1").Value = "Loại thẻ"

# The existing "kikaho" row gains a card-type value.
$ws.Range("C2").Value = "Tiền mặt"

# New account rows.
Set-TextValue $ws.Range("A3") "abcd"
Set-TextValue $ws.Range("B3") "15174"
$ws.Range("C3").Value = "Thẻ ngân hàng"

Set-TextValue $ws.Range("A4") "abcd3"
Set-TextValue $ws.Range("B4") "15174"
$ws.Range("C4").Value = "Tiền mặt"

Set-TextValue $ws.Range("A5") "cxz"
Set-TextValue $ws.Range("B5") "423"
$ws.Range("C5").Value = "Thẻ ngân hàng"
